# Auto-generated edit script: updates Leve profit figures and
# removes stale market-price columns (H:N) for several rows,
# mirroring a scheduled market-data refresh across sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 1205.6757
$ws.Range("I98").Value = 781.38464
$ws.Range("J98").Value = 2208.5454
$ws.Range("K98").Value = 781.38464
$ws.Range("L98").Value = 2208.5454
$ws.Range("M98").Value = 716.61536
$ws.Range("N98").Value = -5204.5454
$ws.Range("H100").Value = 2376.7646
$ws.Range("I100").Value = 1850.625
$ws.Range("J100").Value = 2844.4443
$ws.Range("K100").Value = 1850.625
$ws.Range("L100").Value = 2844.4443
$ws.Range("M100").Value = -1309.625
$ws.Range("N100").Value = -3926.4443
$ws.Range("H122").Value = 1205.6757
$ws.Range("I122").Value = 781.38464
$ws.Range("J122").Value = 2208.5454
$ws.Range("K122").Value = 2344.15392
$ws.Range("L122").Value = 6625.6362
$ws.Range("M122").Value = 105.8460800000003
$ws.Range("N122").Value = -11525.6362
$ws.Range("H125").Value = 984.4545
$ws.Range("I125").Value = 971.5
$ws.Range("K125").Value = 8743.5
$ws.Range("M125").Value = -6283.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 48526.855
$ws.Range("I74").Value = 56402.777
$ws.Range("J74").Value = 1271.3334
$ws.Range("K74").Value = 56402.777
$ws.Range("L74").Value = 1271.3334
$ws.Range("M74").Value = -55528.777
$ws.Range("N74").Value = -3019.3334
$ws.Range("H77").Value = 48526.855
$ws.Range("I77").Value = 56402.777
$ws.Range("J77").Value = 1271.3334
$ws.Range("K77").Value = 282013.885
$ws.Range("L77").Value = 6356.666999999999
$ws.Range("M77").Value = -277645.885
$ws.Range("N77").Value = -15092.667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117:L117").ClearContents()
$ws.Range("N117").ClearContents()
$ws.Range("H118:L118").ClearContents()
$ws.Range("N118").ClearContents()
$ws.Range("H119:L119").ClearContents()
$ws.Range("N119").ClearContents()
$ws.Range("H120:L120").ClearContents()
$ws.Range("N120").ClearContents()
$ws.Range("H122:L122").ClearContents()
$ws.Range("N122").ClearContents()
$ws.Range("H123:L123").ClearContents()
$ws.Range("N123").ClearContents()
$ws.Range("H124:L124").ClearContents()
$ws.Range("N124").ClearContents()
$ws.Range("H125:L125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H126:L126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H127:L127").ClearContents()
$ws.Range("H128:M128").ClearContents()
$ws.Range("H129:L129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H130:L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131:L131").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("H132:L132").ClearContents()
$ws.Range("N132").ClearContents()
$ws.Range("H133:L133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H135:L135").ClearContents()
$ws.Range("H137:L137").ClearContents()
$ws.Range("N137").ClearContents()
$ws.Range("H138:L138").ClearContents()
$ws.Range("N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:L140").ClearContents()
$ws.Range("H141:L141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17858516
$ws.Range("I31").Value = 27027848
$ws.Range("J31").Value = 2446.7896
$ws.Range("K31").Value = 27027848
$ws.Range("L31").Value = 2446.7896
$ws.Range("M31").Value = -27027553
$ws.Range("N31").Value = -3036.7896
$ws.Range("H34").Value = 17858516
$ws.Range("I34").Value = 27027848
$ws.Range("J34").Value = 2446.7896
$ws.Range("K34").Value = 27027848
$ws.Range("L34").Value = 2446.7896
$ws.Range("M34").Value = -27027646
$ws.Range("N34").Value = -2850.7896
$ws.Range("H99").Value = 4363.16
$ws.Range("I99").Value = 3977.8823
$ws.Range("J99").Value = 5181.875
$ws.Range("K99").Value = 3977.8823
$ws.Range("L99").Value = 5181.875
$ws.Range("M99").Value = -2479.8823
$ws.Range("N99").Value = -8177.875
$ws.Range("H126").Value = 4363.16
$ws.Range("I126").Value = 3977.8823
$ws.Range("J126").Value = 5181.875
$ws.Range("K126").Value = 11933.6469
$ws.Range("L126").Value = 15545.625
$ws.Range("M126").Value = -9463.6469
$ws.Range("N126").Value = -20485.625

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 782.55554
$ws.Range("I5").Value = 541.2308
$ws.Range("J5").Value = 1006.6429
$ws.Range("K5").Value = 1623.6924
$ws.Range("L5").Value = 3019.9287
$ws.Range("M5").Value = -1511.6924
$ws.Range("N5").Value = -3243.9287
$ws.Range("H92").Value = 909937.06
$ws.Range("I92").Value = 751.5
$ws.Range("K92").Value = 2254.5
$ws.Range("M92").Value = -1006.5
$ws.Range("H120:M120").ClearContents()
$ws.Range("H121:N121").ClearContents()
$ws.Range("H122:N122").ClearContents()
$ws.Range("H123:N123").ClearContents()
$ws.Range("H124:M124").ClearContents()
$ws.Range("H125:N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:L127").ClearContents()
$ws.Range("N127").ClearContents()
$ws.Range("H128:M128").ClearContents()
$ws.Range("H129:N129").ClearContents()
$ws.Range("H130:N130").ClearContents()
$ws.Range("H131:N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:N133").ClearContents()
$ws.Range("H134:N134").ClearContents()
$ws.Range("H135").Value = 782.55554
$ws.Range("I135").Value = 541.2308
$ws.Range("J135").Value = 1006.6429
$ws.Range("K135").Value = 4871.077200000001
$ws.Range("L135").Value = 9059.786100000001
$ws.Range("M135").Value = -2336.077200000001
$ws.Range("N135").Value = -14129.7861
$ws.Range("H136:N136").ClearContents()
$ws.Range("H137:N137").ClearContents()
$ws.Range("H138:N138").ClearContents()
$ws.Range("H139:N139").ClearContents()
$ws.Range("H140:N140").ClearContents()
$ws.Range("H141:M141").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H125:L125").ClearContents()
$ws.Range("N125").ClearContents()
$ws.Range("H126:N126").ClearContents()
$ws.Range("H127:L127").ClearContents()
$ws.Range("H128:L128").ClearContents()
$ws.Range("N128").ClearContents()
$ws.Range("H129:L129").ClearContents()
$ws.Range("N129").ClearContents()
$ws.Range("H130:L130").ClearContents()
$ws.Range("N130").ClearContents()
$ws.Range("H131:L131").ClearContents()
$ws.Range("N131").ClearContents()
$ws.Range("H132:N132").ClearContents()
$ws.Range("H133:L133").ClearContents()
$ws.Range("N133").ClearContents()
$ws.Range("H134:L134").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("H135:L135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("H136:L136").ClearContents()
$ws.Range("N136").ClearContents()
$ws.Range("H137:L137").ClearContents()
$ws.Range("H138:L138").ClearContents()
$ws.Range("H139:L139").ClearContents()
$ws.Range("N139").ClearContents()
$ws.Range("H140:L140").ClearContents()
$ws.Range("N140").ClearContents()
$ws.Range("H141:L141").ClearContents()
$ws.Range("N141").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 90909640
$ws.Range("I126").Value = 111111570
$ws.Range("J126").Value = 950
$ws.Range("K126").Value = 333334710
$ws.Range("L126").Value = 2850
$ws.Range("M126").Value = -333332240
$ws.Range("N126").Value = -7790

